$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Класс ABC" / analyse_class column (column I) ---------
# Copy the existing header/body cell formatting from column H (Сумма)
# into column I so the new cells reuse the same cell styles, then
# overwrite the copied values with the new header/body text.
$ws.Range("H3").Copy($ws.Range("I3"))
$ws.Range("I3").Value = "Класс ABC"

$ws.Range("H4").Copy($ws.Range("I4"))
$ws.Range("I4").Value = "{`$v->rows[]->analyse_class}"

# New column I is a hidden "helper" column (carries the raw ABC class
# used by conditional formatting/scripting elsewhere), so hide it and
# size it.
$ws.Columns("I").ColumnWidth = 11.66
$ws.Columns("I").Hidden = $true

# --- Window / selection bookkeeping, as left by the editing session ----
$ws.Range("L12").Select()

$win = $excel.ActiveWindow
$win.Left = 5580
